$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BVTs")
$ws.Activate()

# --- Row 8: "Legend Interactivity" ---
# Write values in the same order the original author's Excel session would
# have produced them (this controls the sharedStrings insertion order).
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Legend Interactivity"
$ws.Range("D8").Value = "Click on Legend value."
$ws.Range("E8").Value = "Data should get filtered on clicking particular legend value. Also opacity should get changed for the circle accordingly."
$ws.Range("C8").Value = "Check whether Legend is interactive."

# --- Row 9: "Visual Interactivity" ---
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Visual Interactivity"
$ws.Range("C9").Value = "Check whether Visual is interactive."
$ws.Range("E9").Value = "Data should get filtered for the ""Maths"" category. Also opacity should get changed for that circle accordingly."
$ws.Range("D9").Value = "Click on ""Maths""  path element(circle).                          [NOTE : For overlapping(Intersection) path elements there is no interactivity]                               "

# --- Formatting: copy formats from the rows above which carry the same
# style pattern used for these new rows (A/B = s5, C/D = s4, E = s11) ---
$ws.Range("A7:D7").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A7:D7").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E5").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row heights ---
$ws.Rows.Item(8).RowHeight = 135.75
$ws.Rows.Item(9).RowHeight = 105

# --- Selection, as left by the editing session ---
$ws.Range("D9").Select()
